# Reorder & re-translate the import-parent template header row
# (commit: "improve: ubah urutan dan bahasa header excel parent").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header order / Indonesian wording.
$ws.Range("A1").Value = "Nama Lengkap*"
$ws.Range("B1").Value = "Username*"
$ws.Range("C1").Value = "Pekerjaan"
$ws.Range("D1").Value = "Nomor Telepon"
$ws.Range("E1").Value = "Email"
$ws.Range("F1").Value = "Alamat"
$ws.Range("G1").Value = "Password*"

# H1 used to hold "job"; it's now an empty cell, kept bold & centered like
# the rest of the header row but without the yellow fill / border.
$ws.Range("H1").ClearContents()
$ws.Range("H1").ClearFormats()
$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").HorizontalAlignment = -4108  # xlCenter
$ws.Range("H1").Interior.Pattern = -4142     # xlPatternNone
$ws.Range("H1").Borders.LineStyle = -4142    # xlLineStyleNone

# Move active selection to C4, matching the saved view state.
$ws.Range("C4").Select()
